$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Translate Publisher and Contact values from German to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Fill in the Description value, which was previously empty
$ws.Range("B12").Value = "ID matching types used in consent management "
